$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Inline formula" paragraph: wrap the existing m:oMath in a _GoBack
#    bookmark and replace the trivial "a-b=not c" formula with the full
#    quadratic formula  x = (-b (+-) sqrt(b^2 - 4ac)) / 2a
# ---------------------------------------------------------------------------
$target1 = $d.Content.Find.Execute("Inline formulas are created using the equation editor of Microsoft Word: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$p1 = $d.Paragraphs(5).Range
$xml1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="3954A5B6" w14:textId="77777777" w:rsidR="00424C9C" w:rsidRDefault="00424C9C" w:rsidP="00424C9C"><w:pPr><w:rPr><w:rFonts w:eastAsiaTheme="minorEastAsia"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Inline</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>formulas</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> are </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>created</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>using</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>equation</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> editor of Microsoft Word: </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><m:oMath><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/></w:rPr><m:t>x</m:t></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/></w:rPr><m:t>=</m:t></m:r><m:f><m:fPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr></m:ctrlPr></m:fPr><m:num><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/></w:rPr><m:t>-</m:t></m:r><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/></w:rPr><m:t>b</m:t></m:r><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/></w:rPr><m:t>&#177;</m:t></m:r><m:rad><m:radPr><m:degHide m:val="1"/><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr></m:ctrlPr></m:radPr><m:deg/><m:e><m:sSup><m:sSupPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr></m:ctrlPr></m:sSupPr><m:e><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/></w:rPr><m:t>b</m:t></m:r></m:e><m:sup><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/></w:rPr><m:t>2</m:t></m:r></m:sup></m:sSup><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/></w:rPr><m:t>-4</m:t></m:r><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/></w:rPr><m:t>ac</m:t></m:r></m:e></m:rad></m:num><m:den><m:r><m:rPr><m:sty m:val="p"/></m:rPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/></w:rPr><m:t>2</m:t></m:r><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math" w:cs="Cambria Math"/></w:rPr><m:t>a</m:t></m:r></m:den></m:f></m:oMath><w:bookmarkEnd w:id="0"/></w:p>'
$p1.InsertXML($xml1)

# ---------------------------------------------------------------------------
# 2) "Equazione N: This is an equation" caption paragraph:
#    - drop the centre justification (w:jc)
#    - collapse the spell-checked "This is an equation" runs into a single
#      run (with noProof, since it now sits right after a noProof field run)
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs(10).Range
$xml2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="7858232B" w14:textId="51008E61" w:rsidR="00DA4C09" w:rsidRPr="00DA4C09" w:rsidRDefault="00DA4C09" w:rsidP="00DA4C09"><w:pPr><w:pStyle w:val="Didascalia"/></w:pPr><w:r><w:t xml:space="preserve">Equazione </w:t></w:r><w:fldSimple w:instr=" SEQ Equazione \* ARABIC "><w:r><w:rPr><w:noProof/></w:rPr><w:t>1</w:t></w:r></w:fldSimple><w:r><w:rPr><w:noProof/></w:rPr><w:t xml:space="preserve"> This is an equation</w:t></w:r></w:p>'
$p2.InsertXML($xml2)

# ---------------------------------------------------------------------------
# 3) Drop the trailing empty paragraph and the paragraph that only held the
#    _GoBack bookmark (the bookmark now lives on paragraph 1's formula run).
# ---------------------------------------------------------------------------
$last = $d.Paragraphs($d.Paragraphs.Count).Range
$secondLast = $d.Paragraphs($d.Paragraphs.Count - 1).Range
$delRange = $d.Range($secondLast.Start, $last.End)
$delRange.Delete()
